$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.496595978736877
$ws.Range("B1").Value = 5.209104537963867
$ws.Range("C1").Value = 1.33355176448822
$ws.Range("D1").Value = 0.9230585694313049
$ws.Range("E1").Value = 0.3779242932796478
